$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A. This shifts the existing
# Code / Description / Definition columns one place to the right
# (Code -> B, Description -> C, Definition -> D) and leaves the new
# column A empty, ready to hold the "Version" field.
$ws.Columns.Item(1).Insert()

# New header row: Version | Code | Description | Definition
$ws.Cells.Item(1, 1).Value = "Version"

# Find the last used data row (header is row 1, data starts row 2).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 2).End(-4162).Row

# Fill the new "Version" column with "1.0.0" for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = "1.0.0"
}

$wb.Save()
